$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-22 16:49:21"
$wsZhCn.Range("H4").Value = "2016-03-22 16:49:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-22 16:49:25"
$wsDeDe.Range("H4").Value = "2016-03-22 16:49:53"
